$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 520, shifting existing rows 520:568 down to 521:569
$ws.Rows.Item(520).Insert()

# Populate the new row 520 with data (same constant columns as surrounding rows,
# plus the new varying values from the diff)
$ws.Cells.Item(520, 1).Value = 5
$ws.Cells.Item(520, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(520, 3).Value = "Maule"
$ws.Cells.Item(520, 4).Value = 45132
$ws.Cells.Item(520, 5).Value = 7
$ws.Cells.Item(520, 6).Value = 100114013
$ws.Cells.Item(520, 7).Value = "Zanahoria"
$ws.Cells.Item(520, 8).Value = "Sin especificar"
$ws.Cells.Item(520, 9).Value = "Primera"
$ws.Cells.Item(520, 10).Value = 600
$ws.Cells.Item(520, 11).Value = 5000
$ws.Cells.Item(520, 12).Value = 5000
$ws.Cells.Item(520, 13).Value = 5000
$ws.Cells.Item(520, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(520, 15).Value = "Región de Ñuble"
$ws.Cells.Item(520, 16).Value = 250
$ws.Cells.Item(520, 17).Value = 20
$ws.Cells.Item(520, 18).Value = "Hortaliza"
